# Add 'hole_id' index column to the 'train' sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# Header for the new index column (match the bold/bordered/centered style
# already used by the other header cells in row 1).
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# hole_id values, in row order (rows 2-49), replacing the old 0..47 index.
$holeIds = @(
    "LBU_05_23",
    "LBU_05_05",
    "LBU_05_06",
    "LBU_05_09",
    "MHZ_12_01",
    "LBU_05_25",
    "LBU_98_7",
    "LBU_05_13",
    "MHZ_12_03",
    "LBU_87_5",
    "MHZ_08_03",
    "LBU_05_08",
    "LBU_05_10",
    "MHZ_08_04",
    "LBU_05_12",
    "LBU_01_3",
    "LBU_05_21",
    "LBU_05_15",
    "MHZ_12_02",
    "LBU_05_02",
    "LBU_87_4",
    "LBU_05_20",
    "LBU_05_01",
    "MHZ_08_01",
    "LBU_05_29",
    "LBU_02_4",
    "LBU_87_1",
    "LBU_02_3",
    "LBU_87_2",
    "LBU_07_03",
    "LBU_05_16",
    "LBU_05_22",
    "LBU_98_1",
    "LBU_05_17",
    "MHZ_08_02",
    "LBU_05_27",
    "LBU_98_2",
    "LBU_05_30",
    "LBU_05_03",
    "MHZ_12_04",
    "LBU_05_07",
    "LBU_05_14",
    "LBU_87_6",
    "LBU_05_18",
    "LBU_05_04",
    "LBU_96_2",
    "LBU_05_11",
    "LBU_87_3"
)

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
